$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing values (rows 2-3, 5) ---
$ws.Range("C2").Value = 10.5
$ws.Range("B3").Value = 4
$ws.Range("C3").Value = 9
$ws.Range("C5").Value = 25

# --- Append a new, blank but styled row 6 (same formatting as row 5) ---
$ws.Range("A5:C5").Copy()
$ws.Range("A6:C6").PasteSpecial(-4122)
$ws.Range("A6:C6").ClearContents()

# --- Widen column C slightly to fit the new values (best effort; engine
#     quantizes column widths to 1/7-character pixel steps) ---
$ws.Range("C1").ColumnWidth = 5.15

# --- Move the active selection ---
$ws.Range("E9").Select()
